$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors used for the fill highlights (RGB -> decimal, R + G*256 + B*65536)
$Green  = 45136     # FF00B050 - already used elsewhere (lights)
$Red    = 16711680  # FFFF0000
$Yellow = 16776960  # FFFFFF00
$Purple = 7352480   # FF7030A0

# ---------------------------------------------------------------------
# Existing rows (2 and 3): a handful of cells get new fill highlights.
# Row/col values themselves are unchanged.
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 5).Interior.Color = $Red      # E2 Asena
$ws.Cells.Item(2, 7).Interior.Color = $Yellow   # G2 Jennifer

$ws.Cells.Item(3, 2).Interior.Color = $Red      # B3 Jennifer
$ws.Cells.Item(3, 3).Interior.Color = $Purple   # C3 Dana
$ws.Cells.Item(3, 5).Interior.Color = $Yellow   # E3 Nadja
$ws.Cells.Item(3, 10).Interior.Color = $Yellow  # J3 Laura L.

# ---------------------------------------------------------------------
# Row 4 (new)
# ---------------------------------------------------------------------
$ws.Cells.Item(4, 2).Value = "Anastasia"
$ws.Cells.Item(4, 2).Interior.Color = $Yellow

$ws.Cells.Item(4, 3).Value = "Nadja"
$ws.Cells.Item(4, 3).Interior.Color = $Red

$ws.Cells.Item(4, 4).Value = "Emmy"
$ws.Cells.Item(4, 4).Interior.Color = $Green

$ws.Cells.Item(4, 5).Value = "Jennifer"

$ws.Cells.Item(4, 6).Value = "Gabriela"
$ws.Cells.Item(4, 6).Interior.Color = $Yellow

$ws.Cells.Item(4, 7).Value = "Linda"

$ws.Cells.Item(4, 8).Value = "Asena"

$ws.Cells.Item(4, 9).Value = "Laura L."

$ws.Cells.Item(4, 10).Value = "Laura M."
$ws.Cells.Item(4, 10).Interior.Color = $Red

$ws.Cells.Item(4, 11).Value = "Dana"

$ws.Cells.Item(4, 12).Value = 3

# ---------------------------------------------------------------------
# Row 5 (new) - everyone gets a red highlight except D5 which stays green
# ---------------------------------------------------------------------
$ws.Cells.Item(5, 2).Value = "Jennifer"
$ws.Cells.Item(5, 2).Interior.Color = $Red

$ws.Cells.Item(5, 3).Value = "Nadja"
$ws.Cells.Item(5, 3).Interior.Color = $Red

$ws.Cells.Item(5, 4).Value = "Emmy"
$ws.Cells.Item(5, 4).Interior.Color = $Green

$ws.Cells.Item(5, 5).Value = "Asena"
$ws.Cells.Item(5, 5).Interior.Color = $Red

$ws.Cells.Item(5, 6).Value = "Tara"
$ws.Cells.Item(5, 6).Interior.Color = $Red

$ws.Cells.Item(5, 7).Value = "Gabriela"
$ws.Cells.Item(5, 7).Interior.Color = $Red

$ws.Cells.Item(5, 8).Value = "Linda"
$ws.Cells.Item(5, 8).Interior.Color = $Red

$ws.Cells.Item(5, 9).Value = "Anastasia"
$ws.Cells.Item(5, 9).Interior.Color = $Red

$ws.Cells.Item(5, 10).Value = "Laura M."
$ws.Cells.Item(5, 10).Interior.Color = $Red

$ws.Cells.Item(5, 11).Value = "Laura L."
$ws.Cells.Item(5, 11).Interior.Color = $Red

$ws.Cells.Item(5, 12).Value = 1

# ---------------------------------------------------------------------
# Selection / view bookkeeping to match the saved state
# ---------------------------------------------------------------------
$ws.Range("C3").Select()
